$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The paragraph currently reads "Version 1." and needs to become
# "Version 2." while also changing the internal run layout so that:
#   - "Version" splits into "Versi" + "on" runs
#   - the digit "1" becomes "2"
#   - the trailing "." becomes its own run sitting after the
#     "_GoBack" bookmark instead of being glued onto the " 1" run.
# ------------------------------------------------------------------

# Locate "Version 1." (or an already-updated "Version 2.") so the
# script is resilient to being re-run / to minor offset differences.
$target = $d.Content
$found = $target.Find.Execute("Version 1.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    $target = $d.Content
    $found = $target.Find.Execute("Version 2.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

$paraStart = $target.Start

# Character offsets, relative to the start of the found text:
#   V  e  r  s  i  o  n     1  .
#   0  1  2  3  4  5  6  7  8  9  10
$splitAt   = $paraStart + 5    # between "Versi" and "on"
$digitPos  = $paraStart + 8    # the "1"
$bmNewPos  = $paraStart + 9    # just before the final "."
$dotPos    = $paraStart + 9    # the final "."

# 1. Force "Version" to split into two runs ("Versi" | "on") by
#    briefly dropping a bookmark at the split point, then removing
#    it again -- the bookmark forces a run boundary but leaves no
#    residue once it is deleted.
$splitRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("__tmp_split__", $splitRange)
$d.Bookmarks("__tmp_split__").Delete()

# 2. Change the version digit from "1" to "2".
$digitRange = $d.Range($digitPos, $digitPos + 1)
$digitRange.Text = "2"

# 3. Move the "_GoBack" bookmark so it sits right before the final
#    "." instead of right after it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($bmNewPos, $bmNewPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 4. Re-create the trailing "." as its own run after the bookmark
#    (delete the old one, then insert a fresh run after it).
$dotRange = $d.Range($dotPos, $dotPos + 1)
$dotRange.Delete()
$tailRange = $d.Range($dotPos, $dotPos)
$tailRange.InsertAfter(".")

Write-Output ("Final paragraph text: " + $d.Paragraphs(1).Range.Text)
